$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.277.24"
$ws.Range("E2").Value = "  +1.88%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.53"
$ws.Range("E3").Value = "  +0.60%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.08"
$ws.Range("E5").Value = "  +0.67%  "

# Row 6
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.48%  "

# Row 9
$ws.Range("E9").Value = "  +0.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("E10").Value = "  +1.64%  "

# Row 11
$ws.Range("E11").Value = "  +0.23%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.35"
$ws.Range("E12").Value = "  +0.61%  "

# Row 13
$ws.Range("E13").Value = "  +0.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.634.70"
$ws.Range("E14").Value = "  -0.24%  "

# Row 15
$ws.Range("E15").Value = "  -2.84%  "

# Row 16
$ws.Range("E16").Value = "  +0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.41"
$ws.Range("E17").Value = "  +0.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.272.92"
$ws.Range("E18").Value = "  +1.74%  "

# Row 19
$ws.Range("E19").Value = "  -0.07%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  -0.82%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.84"
$ws.Range("E21").Value = "  +1.58%  "

# Row 22
$ws.Range("E22").Value = "  +1.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.32"
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.60"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  -3.12%  "

# Row 26
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("E27").Value = "  +0.67%  "

# Row 28
$ws.Range("E28").Value = "  +0.48%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("E29").Value = "  +0.93%  "

# Row 30
$ws.Range("E30").Value = "  +0.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0505"
$ws.Range("E31").Value = "  +2.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +0.92%  "

# Row 35
$ws.Range("E35").Value = "  +1.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.914"
$ws.Range("E36").Value = "  +0.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.137.67"
$ws.Range("E37").Value = "  +0.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.555"
$ws.Range("E38").Value = "  +1.74%  "

# Row 40
$ws.Range("E40").Value = "  +1.35%  "

# Row 41
$ws.Range("E41").Value = "  -0.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.60"
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("E43").Value = "  -0.81%  "

# Row 44
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.784.45"
$ws.Range("E45").Value = "  +0.64%  "

# Row 46
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.21"
$ws.Range("E47").Value = "  +3.34%  "

# Row 48
$ws.Range("E48").Value = "  +3.05%  "

# Row 49
$ws.Range("E49").Value = "  +2.81%  "

# Row 50
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.418"
$ws.Range("E50").Value = "  +0.29%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.67"
$ws.Range("E51").Value = "  +2.87%  "
